$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet1 (展览)
$ws1.Range("C2").Value = "苏州·暑假COS动漫展-CF01（取消）"
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F5").Value = 277
$ws1.Range("F6").Value = 1082
$ws1.Range("F7").Value = 1420
$ws1.Range("F8").Value = 584
$ws1.Range("F9").Value = 104
$ws1.Range("F10").Value = 743
$ws1.Range("F11").Value = 69
$ws1.Range("F12").Value = 141
$ws1.Range("C13").Value = "苏州·环球港动漫节（取消）"
$ws1.Range("G13").Value = "不可售"
$ws1.Range("F14").Value = 425
$ws1.Range("F15").Value = 1329
$ws1.Range("F17").Value = 94
$ws1.Range("F18").Value = 272
$ws1.Range("F20").Value = 645
$ws1.Range("F21").Value = 31
$ws1.Range("F22").Value = 207
$ws1.Range("F23").Value = 15
$ws1.Range("F24").Value = 5722
$ws1.Range("F25").Value = 57
$ws1.Range("F26").Value = 119
$ws1.Range("F27").Value = 91
$ws1.Range("F29").Value = 14296
$ws1.Range("F30").Value = 1424
$ws1.Range("F31").Value = 199
$ws1.Range("F32").Value = 95
$ws1.Range("F34").Value = 1391
$ws1.Range("F35").Value = 593
$ws1.Range("F36").Value = 4183
$ws1.Range("F37").Value = 126
$ws1.Range("F38").Value = 354

# Sheet2 (演出)
$ws2.Range("F2").Value = 342

# Sheet4 (全部类型)
$ws4.Range("C2").Value = "苏州·暑假COS动漫展-CF01（取消）"
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F5").Value = 277
$ws4.Range("F6").Value = 1082
$ws4.Range("F7").Value = 1420
$ws4.Range("F8").Value = 584
$ws4.Range("F9").Value = 104
$ws4.Range("F10").Value = 743
$ws4.Range("F11").Value = 69
$ws4.Range("F12").Value = 141
$ws4.Range("C13").Value = "苏州·环球港动漫节（取消）"
$ws4.Range("G13").Value = "不可售"
$ws4.Range("F14").Value = 425
$ws4.Range("F15").Value = 1329
$ws4.Range("F17").Value = 94
$ws4.Range("F18").Value = 272
$ws4.Range("F19").Value = 342
$ws4.Range("F21").Value = 645
$ws4.Range("F23").Value = 31
$ws4.Range("F24").Value = 207
$ws4.Range("F25").Value = 15
$ws4.Range("F27").Value = 5722
$ws4.Range("F28").Value = 57
$ws4.Range("F29").Value = 119
$ws4.Range("F30").Value = 91
$ws4.Range("F32").Value = 14296
$ws4.Range("F33").Value = 1424
$ws4.Range("F34").Value = 199
$ws4.Range("F35").Value = 95
$ws4.Range("F37").Value = 1391
$ws4.Range("F38").Value = 593
$ws4.Range("F39").Value = 4183
$ws4.Range("F40").Value = 126
$ws4.Range("F41").Value = 354

